$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("E3").Value = "Callback function that runs once the animation has fully completed. Does not run if animation is replaced by another ``goPhase()`` call."
$ws.Range("D2").Value = "[CodePen](https://codepen.io/maiCoding/pen/JexvLW)"
$ws.Range("D3").Value = "[CodePen](https://codepen.io/maiCoding/pen/JexvLW)"

$ws.Range("E16").Select()
$excel.ActiveWindow.ScrollRow = 1
$excel.ActiveWindow.ScrollColumn = 2
